$wb = $excel.ActiveWorkbook

# --- Worksheet 7 "text_coercion": add column B explanatory strings ---
$ws7 = $wb.Worksheets.Item(7)

$ws7.Range("B1").Value = "explanation"
$ws7.Range("B2").Value = "text"
$ws7.Range("B3").Value = "blank"
$ws7.Range("B4").Value = "logical F"
$ws7.Range("B5").Value = "boolean"
$ws7.Range("B6").Value = "floating point"
$ws7.Range("B7").Value = "date"
$ws7.Range("B8").Value = "text"

# Move the active cell / selection to B9 and make this sheet the active (tab-selected) one,
# mirroring the selection state saved in the edited workbook.
[void]$ws7.Range("B9").Select()

# --- Workbook-level view state: active tab now points at text_coercion (index 6, 0-based) ---
$win = $excel.Windows.Item(1)
$win.Left = 2880
$win.Top = 1160
